$d = $word.ActiveDocument

# --- 1) Title: "Практическая работа 2" -> "Практическая работа 4",
#     plus a collapsed _GoBack bookmark placed right after that run.
$d.Content.Find.Execute("Практическая работа 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Практическая работа 4", 2) | Out-Null

$titleRange = $d.Content
$titleRange.Find.Execute("Практическая работа 4") | Out-Null
$titleEnd = $titleRange.End

# A *collapsed* bookmark sitting exactly at the last valid offset of a
# paragraph's text (one before the paragraph mark) gets mis-anchored to the
# paragraph's start by this host. Work around it: insert a throwaway
# character right after the run so the insertion point is no longer the
# paragraph's last position, drop the bookmark there, then delete the
# throwaway character again.
$d.Range($titleEnd, $titleEnd).InsertAfter("X") | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($titleEnd, $titleEnd)) | Out-Null
$d.Range($titleEnd, $titleEnd + 1).Delete() | Out-Null

# --- 2) Body paragraph: merge
#     "написал код для просмотра данных, фильтрации" + " по коду и жанру" +
#     ", а также код для сброса фильтров" (with a _GoBack bookmark sitting
#     between the 2nd and 3rd run) into one run / one <w:t>.
$r1 = $d.Content
$r1.Find.Execute("написал код для просмотра данных, фильтрации") | Out-Null
$runEnd = $r1.End

$r2 = $d.Content
$r2.Find.Execute(", а также код для сброса фильтров") | Out-Null
$tailEnd = $r2.End

# Same end-of-paragraph-minus-one quirk as above can bite Range.Text
# assignment; nudge the end forward to the real paragraph end when so.
$para = $r2.Paragraphs(1)
if ($tailEnd -eq ($para.Range.End - 1)) {
    $tailEnd = $para.Range.End
}

$full = $d.Range($runEnd, $tailEnd)
$full.Text = " по коду и жанру, а также код для сброса фильтров"
